# Automatic post-commit hook for streamlit
#
# The analysis dataset had two records (sampleid 280 in row 7 and
# sampleid 424 in row 8) whose data got swapped between the two rows.
# Reproduce that by swapping the entire row contents of row 7 and row 8,
# using Copy() (instead of literal re-typing of values) so that text vs.
# number cell typing, number formats and styles are preserved exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowA = $ws.Range("A7:CP7")
$rowB = $ws.Range("A8:CP8")
$scratch = $ws.Range("A1000:CP1000")

$rowA.Copy($scratch)
$rowB.Copy($rowA)
$scratch.Copy($rowB)
$scratch.Clear()
